$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 75 with revised figures ---
$ws.Cells.Item(75, 4).Value = 1977
$ws.Cells.Item(75, 5).Value = 3072
$ws.Cells.Item(75, 19).Value = 67442
$ws.Cells.Item(75, 20).Value = 1961
$ws.Cells.Item(75, 21).Value = 17181
$ws.Cells.Item(75, 22).Value = 73459

# --- Append new row 76 for period 01-04-2021 ---
# Force text formatting first so Excel does not auto-convert the
# dd-mm-yyyy-looking label into a date serial number, then restore the
# cell to the workbook's default (unstyled) look.
$ws.Cells.Item(76, 1).NumberFormat = "@"
$ws.Cells.Item(76, 1).Value = "01-04-2021"
$ws.Cells.Item(76, 1).Style = "Normal"
$ws.Cells.Item(76, 2).Value = 148520
$ws.Cells.Item(76, 3).Value = 4590
$ws.Cells.Item(76, 4).Value = 2146
$ws.Cells.Item(76, 5).Value = 2444
$ws.Cells.Item(76, 6).Value = 62410
$ws.Cells.Item(76, 7).Value = 3314
$ws.Cells.Item(76, 8).Value = 59096
$ws.Cells.Item(76, 9).Value = 98
$ws.Cells.Item(76, 10).Value = 81422
$ws.Cells.Item(76, 11).Value = 18009
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = 63413
$ws.Cells.Item(76, 14).Value = 0
$ws.Cells.Item(76, 15).Value = 148520
$ws.Cells.Item(76, 16).Value = 147978
$ws.Cells.Item(76, 17).Value = 542
$ws.Cells.Item(76, 18).Value = 148520
$ws.Cells.Item(76, 19).Value = 66707
$ws.Cells.Item(76, 20).Value = 1893
$ws.Cells.Item(76, 21).Value = 14714
$ws.Cells.Item(76, 22).Value = 65206
